# Updates the cryptos worksheet with refreshed price / volume data,
# matching the author's commit "Updated cryptos list ... with GitHub Actions".
# Also accounts for 3 coins (Polkadot / Chainlink / WrappedEther) having
# swapped rank positions (rows 13-15) between the two snapshots.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text, avoiding Excel's automatic
# number/date coercion for numeric-looking strings (e.g. "0.4598",
# "1.001"), while leaving the cell's style untouched (no stray
# number-format / style index left behind on the cell).
function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$Text
    )
    $range = $ws.Range($CellRef)
    $range.NumberFormat = "@"
    $range.Value = $Text
    $range.ClearFormats()
}

Set-TextValue "D2" '28.293.91'
Set-TextValue "E2" '  -5.61%  '
Set-TextValue "D3" '1.835.72'
Set-TextValue "E3" '  -5.44%  '
Set-TextValue "E4" '  -0.75%  '
Set-TextValue "D5" '331.10'
Set-TextValue "E5" '  -1.23%  '
Set-TextValue "E6" '  -0.69%  '
Set-TextValue "D7" '0.4598'
Set-TextValue "E7" '  -5.03%  '
Set-TextValue "D8" '0.3862'
Set-TextValue "E8" '  -6.73%  '
Set-TextValue "D9" '45.81'
Set-TextValue "E9" '  -3.79%  '
Set-TextValue "D10" '0.07854'
Set-TextValue "E10" '  -4.30%  '
Set-TextValue "D11" '0.9663'
Set-TextValue "E11" '  -4.98%  '
Set-TextValue "D12" '21.87'
Set-TextValue "E12" '  -8.51%  '
Set-TextValue "B13" 'WrappedEther'
Set-TextValue "C13" 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue "D13" '1.823.39'
Set-TextValue "E13" '  -7.22%  '
Set-TextValue "B14" 'Polkadot'
Set-TextValue "C14" 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue "D14" '5.727'
Set-TextValue "E14" '  -6.23%  '
Set-TextValue "B15" 'Chainlink'
Set-TextValue "C15" 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue "D15" '6.926'
Set-TextValue "E15" '  -5.37%  '
Set-TextValue "D16" '0.06865'
Set-TextValue "E16" '  -0.05%  '
Set-TextValue "D17" '1.001'
Set-TextValue "D18" '86.70'
Set-TextValue "E18" '  -5.17%  '
Set-TextValue "D19" '0.000009919'
Set-TextValue "E19" '  -4.65%  '
Set-TextValue "D20" '16.90'
Set-TextValue "E20" '  -5.48%  '
Set-TextValue "D22" '28.296.92'
Set-TextValue "E22" '  -5.59%  '
Set-TextValue "D23" '5.339'
Set-TextValue "E23" '  -5.52%  '
Set-TextValue "D24" '11.00'
Set-TextValue "E24" '  -7.73%  '
Set-TextValue "D25" '2.155'
Set-TextValue "E25" '  -1.69%  '
Set-TextValue "D26" '2.010.71'
Set-TextValue "E26" '  -8.59%  '
Set-TextValue "D27" '153.14'
Set-TextValue "E27" '  -2.52%  '
Set-TextValue "D28" '19.22'
Set-TextValue "E28" '  -4.59%  '
Set-TextValue "D29" '5.800'
Set-TextValue "E29" '  -13.13%  '
Set-TextValue "D30" '1.976'
Set-TextValue "E30" '  -6.21%  '
Set-TextValue "D31" '116.75'
Set-TextValue "E31" '  -3.86%  '
Set-TextValue "D32" '0.9441'
Set-TextValue "E32" '  -7.23%  '
Set-TextValue "E33" '  -3.36%  '
Set-TextValue "D34" '5.291'
Set-TextValue "E34" '  -5.93%  '
Set-TextValue "D35" '3.445'
Set-TextValue "E35" '  -3.14%  '
Set-TextValue "D36" '1.326'
Set-TextValue "E36" '  -6.92%  '
Set-TextValue "D37" '0.06036'
Set-TextValue "E37" '  -7.97%  '
Set-TextValue "D38" '0.02155'
Set-TextValue "E38" '  -6.10%  '
Set-TextValue "E39" '  -5.50%  '
Set-TextValue "D40" '0.9997'
Set-TextValue "E40" '  -0.83%  '
Set-TextValue "D41" '0.5620'
Set-TextValue "E41" '  -6.11%  '
Set-TextValue "D42" '7.559'
Set-TextValue "E42" '  -5.73%  '
Set-TextValue "D43" '10.01'
Set-TextValue "E43" '  -6.83%  '
Set-TextValue "D44" '0.1782'
Set-TextValue "E44" '  -3.99%  '
Set-TextValue "E45" '  -2.90%  '
Set-TextValue "D46" '2.274'
Set-TextValue "E46" '  -10.36%  '
Set-TextValue "D47" '11.68'
Set-TextValue "E47" '  -6.17%  '
Set-TextValue "D48" '0.5298'
Set-TextValue "E48" '  -5.00%  '
Set-TextValue "D49" '0.07028'
Set-TextValue "E49" '  -6.52%  '
Set-TextValue "D50" '1.841'
Set-TextValue "E50" '  -7.40%  '
Set-TextValue "D51" '113.17'
Set-TextValue "E51" '  -3.92%  '
